$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.778.86"
$ws.Range("E2").Value = "  -7.89%  "

$ws.Range("D3").Value = "2.519.61"
$ws.Range("E3").Value = "  -3.73%  "

$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "296.53"
$ws.Range("E5").Value = "  -3.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.10"
$ws.Range("E6").Value = "  -7.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.569"
$ws.Range("E7").Value = "  -5.61%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.547"
$ws.Range("E9").Value = "  -5.95%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.23"
$ws.Range("E10").Value = "  -8.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0798"
$ws.Range("E11").Value = "  -5.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.63"
$ws.Range("E12").Value = "  -6.33%  "

$ws.Range("E13").Value = "  +0.33%  "

$ws.Range("D14").Value = "2.898.92"
$ws.Range("E14").Value = "  -3.72%  "

$ws.Range("D15").Value = "2.514.37"
$ws.Range("E15").Value = "  -3.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.865"
$ws.Range("E16").Value = "  -6.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.98"
$ws.Range("E17").Value = "  -6.83%  "

$ws.Range("D18").Value = "42.715.89"
$ws.Range("E18").Value = "  -8.27%  "

$ws.Range("D19").Value = "0.0₃0961"
$ws.Range("E19").Value = "  -5.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.51"
$ws.Range("E20").Value = "  -3.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.27"
$ws.Range("E21").Value = "  -5.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.29"
$ws.Range("E22").Value = "  +0.57%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "258.64"
$ws.Range("E23").Value = "  -6.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.88"
$ws.Range("E24").Value = "  -5.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.13"
$ws.Range("E25").Value = "  -2.74%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.01"
$ws.Range("E26").Value = "  -2.09%  "

$ws.Range("E27").Value = "  +0.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.92"
$ws.Range("E28").Value = "  -6.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  -0.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.86"
$ws.Range("E30").Value = "  -4.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.99"
$ws.Range("E31").Value = "  -6.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.46"
$ws.Range("E32").Value = "  -5.00%  "

$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.19"
$ws.Range("E33").Value = "  -2.78%  "

$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.69"
$ws.Range("E34").Value = "  -0.99%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.76"
$ws.Range("E35").Value = "  -2.81%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0795"
$ws.Range("E36").Value = "  -5.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.114"
$ws.Range("E37").Value = "  -7.32%  "

$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.118"
$ws.Range("E38").Value = "  -3.96%  "

$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "23.84"
$ws.Range("E39").Value = "  +1.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.23"
$ws.Range("E40").Value = "  +1.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.42"
$ws.Range("E41").Value = "  -5.83%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0308"
$ws.Range("E42").Value = "  -7.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.82"
$ws.Range("E43").Value = "  -5.59%  "

$ws.Range("D44").Value = "2.012.28"
$ws.Range("E44").Value = "  -5.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.995"
$ws.Range("E45").Value = "  -0.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "85.54"
$ws.Range("E46").Value = "  -9.28%  "

$ws.Range("E47").Value = "  +2.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.86"
$ws.Range("E48").Value = "  -7.21%  "

$ws.Range("D49").Value = "2.749.81"
$ws.Range("E49").Value = "  -4.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.75"
$ws.Range("E50").Value = "  -6.55%  "

$ws.Range("E51").Value = "  -7.82%  "
